# Restructured pathExport, covered 'special' cases/rmTypes and defined the output format
#
# This script edits the "FLAT_Paths" worksheet of the workbook:
#  1. Clears the "Mandatory Paths" column (C) for all data rows (rows 2-29),
#     leaving only the FLAT-Path list (column A) and the header row intact.
#  2. Swaps the ordering of the "|code" / "|terminology" pairs for each
#     language|.. / encoding|.. / territory|.. entry so that "|terminology"
#     precedes "|code" in the dropdown list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLAT_Paths")

# 1. Clear the "Mandatory Paths" column (C) for the data rows (keep header C1).
$ws.Range("C2:C29").ClearContents()

# 2. Swap adjacent row pairs in column A so "|terminology" comes before "|code".
$swapRowPairs = @(
    @(29, 30),
    @(31, 32),
    @(41, 42),
    @(43, 44),
    @(52, 53),
    @(54, 55),
    @(63, 64),
    @(65, 66),
    @(79, 80),
    @(81, 82),
    @(87, 88),
    @(89, 90)
)

foreach ($pair in $swapRowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $cell1 = $ws.Cells.Item($r1, 1)
    $cell2 = $ws.Cells.Item($r2, 1)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
}
